$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 contains a bogus date value (serial 61 = 1900-03-01 boundary bug,
# i.e. a date before 1900-03-01) - delete the whole row so the remaining
# records shift up.
$ws.Rows.Item(2).Delete()

# Reselect row 2 (now holding the former row-3 data) to match the saved
# workbook's selection state.
$ws.Rows.Item(2).Select()
